$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.486150026321411
$ws.Range("B1").Value = 2.608836650848389
$ws.Range("C1").Value = 3.88152027130127
$ws.Range("D1").Value = 3.940397262573242
$ws.Range("E1").Value = 2.746790409088135
